$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE (keep as literal text, matching the original "Text"-style entry)
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "05/03/2021"

# START TIME
$ws.Range("F4").Value = "07:00"

# WEATHER
$ws.Range("B5").Value = "Sunny"

# END TIME
$ws.Range("F5").Value = "16:00"

# Manpower / Equipment table row 8
$ws.Range("A8").Value = "Exbon Development Inc."
$ws.Range("B8").Value = "Carpenter"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = "Forklift"
$ws.Range("F8").Value = "Sneeze Partition Installation"

# Manpower / Equipment table row 9
$ws.Range("A9").Value = "JPUS"
$ws.Range("B9").Value = "Laborer"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = "Dump Truck"
$ws.Range("F9").Value = "Existing Partition Removal"

# TESTS & INSPECTIONS note
$ws.Range("A18").Value = "Inspection is scheduled on May 4, 2021 at 3PM. "

# CORRECTIONAL ITEMS note
$ws.Range("A23").Value = "New frosted panel has a crack. Needs to be replaced."

# NOTE
$ws.Range("A28").Value = "All punchwork need to be performed during off-hours."
